$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits at the end
#    of the "...func_tests." sentence (it will be re-created further
#    below, at the spot of the real last edit).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Insert the missing comma: "...может присутствовать число цифр..."
#    -> "...может присутствовать, число цифр..."
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("может присутствовать", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor text 'может присутствовать' not found"
}
$endPos = $r.End

# Insert the comma right after "присутствовать" (before the existing
# space that separates it from "число").
$insertionPoint = $d.Range($endPos, $endPos)
$insertionPoint.InsertAfter(",")

# The comma plus the (pre-existing) following space now occupy
# [$endPos, $endPos + 2). Force that pair onto its own run (mirroring
# how Word itself breaks runs around freshly-typed text) by briefly
# bookmarking it and dropping the helper bookmark again.
$commaSpaceRange = $d.Range($endPos, $endPos + 2)
if ($d.Bookmarks.Exists("zzTmpSplit")) {
    $d.Bookmarks("zzTmpSplit").Delete()
}
$d.Bookmarks.Add("zzTmpSplit", $commaSpaceRange) | Out-Null
$d.Bookmarks("zzTmpSplit").Delete()

# ------------------------------------------------------------------
# 3) Drop the "_GoBack" bookmark exactly where the edit happened: right
#    after the newly inserted ", " and before "число".
# ------------------------------------------------------------------
$goBackPos = $d.Range($endPos + 2, $endPos + 2)
$d.Bookmarks.Add("_GoBack", $goBackPos) | Out-Null
